$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# --- Row 1 (header row): replace data-duplicate values with proper column
# labels, matching the header row used on the other property sheets
# (name / capacity / owner / register_date / register_reason /
#  acquire_value / property_category / category / date / legislator_name /
#  legislator_id / source_file / index).
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# New header cells need the same bold/bordered style as the existing
# header cells (B1:G1 already carry it).
$ws.Range("B1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)

# --- Row 2 / Row 3: existing car data (name/capacity/owner/register_date/
# register_reason/acquire_value) is untouched; append the metadata columns
# that are present on every other property sheet.
$ws.Range("H2:H3").Value = "land"
$ws.Range("I2:I3").Value = "normal"
# Force text so "2012-04-12" isn't re-interpreted as a date serial.
$ws.Range("J2:J3").NumberFormat = "@"
$ws.Range("J2:J3").Value = "2012-04-12"
$ws.Range("K2:K3").Value = "林正二"
$ws.Range("L2:L3").Value = 788
$ws.Range("M2:M3").Value = "tmp32921"
$ws.Range("N2").Value = 40
$ws.Range("N3").Value = 41

# Match the plain data-row style already used by B2:G3 (also clears the
# temporary text number-format applied above).
$ws.Range("B2").Copy()
$ws.Range("H2:N3").PasteSpecial(-4122)
